$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'64.678.13"
$ws.Range("E2").Value = "  +2.73%  "

# Row 3
$ws.Range("D3").Value = "'3.369.21"
$ws.Range("E3").Value = "  +1.76%  "

# Row 4
$ws.Range("E4").Value = "  -0.22%  "

# Row 5
$ws.Range("D5").Value = "'562.97"
$ws.Range("E5").Value = "  +2.73%  "

# Row 6
$ws.Range("D6").Value = "'176.18"
$ws.Range("E6").Value = "  +4.69%  "

# Row 7
$ws.Range("D7").Value = "'0.630"
$ws.Range("E7").Value = "  +4.70%  "

# Row 8
$ws.Range("D8").Value = "'3.360.17"
$ws.Range("E8").Value = "  +1.82%  "

# Row 9
$ws.Range("E9").Value = "  +0.00%  "

# Row 10
$ws.Range("D10").Value = "'0.634"
$ws.Range("E10").Value = "  +5.41%  "

# Row 11
$ws.Range("E11").Value = "  +11.56%  "

# Row 12
$ws.Range("D12").Value = "'55.08"
$ws.Range("E12").Value = "  +2.69%  "

# Row 13
$ws.Range("D13").Value = "'0.0000276"
$ws.Range("E13").Value = "  +6.13%  "

# Row 14
$ws.Range("E14").Value = "  +4.26%  "

# Row 15
$ws.Range("D15").Value = "'3.902.34"
$ws.Range("E15").Value = "  +0.71%  "

# Row 16
$ws.Range("E16").Value = "  +5.01%  "

# Row 17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.118"
$ws.Range("E17").Value = "  +1.62%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'3.358.78"
$ws.Range("E18").Value = "  -0.70%  "

# Row 19
$ws.Range("D19").Value = "'11.87"
$ws.Range("E19").Value = "  +4.05%  "

# Row 20
$ws.Range("D20").Value = "'64.566.91"
$ws.Range("E20").Value = "  +2.45%  "

# Row 21
$ws.Range("E21").Value = "  +3.82%  "

# Row 22
$ws.Range("D22").Value = "'469.92"
$ws.Range("E22").Value = "  +19.26%  "

# Row 23
$ws.Range("E23").Value = "  +15.13%  "

# Row 24
$ws.Range("D24").Value = "'4.13"
$ws.Range("E24").Value = "  +3.80%  "

# Row 25
$ws.Range("D25").Value = "'86.68"
$ws.Range("E25").Value = "  +7.15%  "

# Row 26
$ws.Range("D26").Value = "'13.48"
$ws.Range("E26").Value = "  +4.36%  "

# Row 27
$ws.Range("E27").Value = "  +1.96%  "

# Row 28
$ws.Range("D28").Value = "'2.85"
$ws.Range("E28").Value = "  +5.85%  "

# Row 29
$ws.Range("D29").Value = "'8.83"
$ws.Range("E29").Value = "  +4.26%  "

# Row 30
$ws.Range("D30").Value = "'30.36"
$ws.Range("E30").Value = "  +5.80%  "

# Row 31
$ws.Range("D31").Value = "'6.66"
$ws.Range("E31").Value = "  +4.87%  "

# Row 32
$ws.Range("D32").Value = "'11.52"
$ws.Range("E32").Value = "  +3.53%  "

# Row 33
$ws.Range("D33").Value = "'579.58"
$ws.Range("E33").Value = "  -0.77%  "

# Row 34
$ws.Range("E34").Value = "  +4.55%  "

# Row 35
$ws.Range("D35").Value = "'59.96"
$ws.Range("E35").Value = "  +3.82%  "

# Row 36
$ws.Range("E36").Value = "  +0.05%  "

# Row 37
$ws.Range("E37").Value = "  -5.41%  "

# Row 38
$ws.Range("D38").Value = "'35.92"
$ws.Range("E38").Value = "  +1.79%  "

# Row 39
$ws.Range("E39").Value = "  +4.43%  "

# Row 40
$ws.Range("E40").Value = "  +0.22%  "

# Row 41
$ws.Range("E41").Value = "  +2.88%  "

# Row 42
$ws.Range("D42").Value = "'3.086.14"
$ws.Range("E42").Value = "  -1.22%  "

# Row 43
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "  -0.47%  "

# Row 44
$ws.Range("D44").Value = "'2.83"
$ws.Range("E44").Value = "  +2.64%  "

# Row 45
$ws.Range("E45").Value = "  +3.54%  "

# Row 46
$ws.Range("E46").Value = "  +4.53%  "

# Row 47
$ws.Range("E47").Value = "  +1.08%  "

# Row 48
$ws.Range("E48").Value = "  +5.22%  "

# Row 49
$ws.Range("D49").Value = "'2.59"
$ws.Range("E49").Value = "  +0.44%  "

# Row 50
$ws.Range("D50").Value = "'137.44"
$ws.Range("E50").Value = "  +4.71%  "

# Row 51
$ws.Range("D51").Value = "'8.41"
$ws.Range("E51").Value = "  +5.84%  "
